$wb = $excel.ActiveWorkbook

# Sheet1 = "번역" (translation list), Sheet2 = "용어통일" (glossary/notes list)
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Add new glossary note "원형 경로" as a new row in the notes sheet first,
#     so the shared-string table ordering matches the target workbook. ---
$ws2.Range("A14").Value = "원형 경로"

# --- Sheet1 corrections: add a space inside two existing Korean translations ---
# Row 12 ("closed curve")
$ws1.Range("B12").Value = "닫힌 곡선"
# Row 59 ("open disk")
$ws1.Range("B59").Value = "열린 원판"

# --- Append two more new glossary notes to Sheet2 ---
$ws2.Range("A15").Value = "단순연결 영역"
$ws2.Range("A16").Value = "단위 원판"

# --- Update the selections to reflect where the editor last left off ---
# Select sheet1's edited cell first...
$ws1.Range("B59").Select() | Out-Null
# ...then select sheet2's new last row + 1 (keeps Sheet2 as the active tab)
$ws2.Range("A17").Select() | Out-Null

Write-Host "Applied word_list.xlsx glossary updates"
